# Generate Report for Handoff
#
# A new handoff round completed for 82752071-1c07-4ffe-91f3-7a0b4128ebe5.
# This moves that file's row to the end of the "ready" block (rows 6-9) on
# every sheet, shifting the other three rows (a5d06b60, 2c19810b, 57e7464f)
# up by one, and refreshes 82752071's handoff timestamps.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: columns A-G
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A6").Value = "a5d06b60-17cf-44f3-b757-6f4b919228e3.md"
$ws.Range("B6").Value = "e2e\a5d06b60-17cf-44f3-b757-6f4b919228e3.md"
$ws.Range("C6").Value = ".md"
$ws.Range("E6").Value = "In Translation"
$ws.Range("F6").Value = "In Translation"
$ws.Range("G6").Value = "2016-08-13 10:50:01"

$ws.Range("A7").Value = "2c19810b-ee51-452b-aafa-8f125547ddd1.md"
$ws.Range("B7").Value = "e2e\2c19810b-ee51-452b-aafa-8f125547ddd1.md"
$ws.Range("C7").Value = ".md"
$ws.Range("E7").Value = "Ready for handoff"
$ws.Range("F7").Value = "Ready for handoff"
$ws.Range("G7").Value = "2016-08-13 10:48:31"

$ws.Range("A8").Value = "57e7464f-7052-4509-b5fa-2d4ffc83dec5.md"
$ws.Range("B8").Value = "e2e\57e7464f-7052-4509-b5fa-2d4ffc83dec5.md"
$ws.Range("C8").Value = ".md"
$ws.Range("E8").Value = "Ready for handoff"
$ws.Range("F8").Value = "Ready for handoff"
$ws.Range("G8").Value = "2016-08-13 10:51:33"

$ws.Range("A9").Value = "82752071-1c07-4ffe-91f3-7a0b4128ebe5.md"
$ws.Range("B9").Value = "e2e\82752071-1c07-4ffe-91f3-7a0b4128ebe5.md"
$ws.Range("C9").Value = ".md"
$ws.Range("E9").Value = "Ready for handoff"
$ws.Range("F9").Value = "Ready for handoff"
$ws.Range("G9").Value = "2016-08-13 10:55:21"

# ---------------------------------------------------------------------
# zh-cn sheet: columns A-H (I-P unchanged/blank for these rows)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A6").Value = "a5d06b60-17cf-44f3-b757-6f4b919228e3.md"
$ws.Range("C6").Value = "In Translation"
$ws.Range("G6").Value = "a5d06b60-17cf-44f3-b757-6f4b919228e3.30e8ead1a3b0dc3f3d31e07b2151d916235d57a2.zh-cn.xlf"
$ws.Range("H6").Value = "2016-08-13 10:49:53"

$ws.Range("A7").Value = "2c19810b-ee51-452b-aafa-8f125547ddd1.md"
$ws.Range("C7").Value = "Ready for handoff"
$ws.Range("G7").Value = "2c19810b-ee51-452b-aafa-8f125547ddd1.b3f5881bb1f0faf01175652c3d302b3e1a37cad7.zh-cn.xlf"
$ws.Range("H7").Value = "2016-08-13 10:48:24"

$ws.Range("A8").Value = "57e7464f-7052-4509-b5fa-2d4ffc83dec5.md"
$ws.Range("C8").Value = "Ready for handoff"
$ws.Range("G8").Value = "57e7464f-7052-4509-b5fa-2d4ffc83dec5.93a72b73dff4fbf0545eafe0775adbb35b50061a.zh-cn.xlf"
$ws.Range("H8").Value = "2016-08-13 10:51:25"

$ws.Range("A9").Value = "82752071-1c07-4ffe-91f3-7a0b4128ebe5.md"
$ws.Range("C9").Value = "Ready for handoff"
$ws.Range("G9").Value = "82752071-1c07-4ffe-91f3-7a0b4128ebe5.6e7fac94290f194893159599289946e9174f74b9.zh-cn.xlf"
$ws.Range("H9").Value = "2016-08-13 10:55:13"

# ---------------------------------------------------------------------
# de-de sheet: columns A-H (I-P unchanged/blank for these rows)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A6").Value = "a5d06b60-17cf-44f3-b757-6f4b919228e3.md"
$ws.Range("C6").Value = "In Translation"
$ws.Range("G6").Value = "a5d06b60-17cf-44f3-b757-6f4b919228e3.30e8ead1a3b0dc3f3d31e07b2151d916235d57a2.de-de.xlf"
$ws.Range("H6").Value = "2016-08-13 10:50:01"

$ws.Range("A7").Value = "2c19810b-ee51-452b-aafa-8f125547ddd1.md"
$ws.Range("C7").Value = "Ready for handoff"
$ws.Range("G7").Value = "2c19810b-ee51-452b-aafa-8f125547ddd1.b3f5881bb1f0faf01175652c3d302b3e1a37cad7.de-de.xlf"
$ws.Range("H7").Value = "2016-08-13 10:48:31"

$ws.Range("A8").Value = "57e7464f-7052-4509-b5fa-2d4ffc83dec5.md"
$ws.Range("C8").Value = "Ready for handoff"
$ws.Range("G8").Value = "57e7464f-7052-4509-b5fa-2d4ffc83dec5.93a72b73dff4fbf0545eafe0775adbb35b50061a.de-de.xlf"
$ws.Range("H8").Value = "2016-08-13 10:51:33"

$ws.Range("A9").Value = "82752071-1c07-4ffe-91f3-7a0b4128ebe5.md"
$ws.Range("C9").Value = "Ready for handoff"
$ws.Range("G9").Value = "82752071-1c07-4ffe-91f3-7a0b4128ebe5.6e7fac94290f194893159599289946e9174f74b9.de-de.xlf"
$ws.Range("H9").Value = "2016-08-13 10:55:21"
